# Weekly update: insert two new report rows (week of 2022-02-07, serial 44610)
# for "Betarraga" at Vega Central Mapocho de Santiago, pushing the existing
# rows 397:419 down to 399:421.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 397; everything that was on
# rows 397:419 shifts down to 399:421 (and picks up formatting, incl. the
# date NumberFormat on column D, from the row above).
$ws.Rows("397:398").Insert()

# --- New row 397 ---
$ws.Range("A397").Value = 9
$ws.Range("B397").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C397").Value = "Metropolitana"
$ws.Range("D397").Value2 = 44610
$ws.Range("E397").Value = 13
$ws.Range("F397").Value = 100114014
$ws.Range("G397").Value = "Betarraga"
$ws.Range("H397").Value = "Sin especificar"
$ws.Range("I397").Value = "Primera"
$ws.Range("J397").Value = 7900
$ws.Range("K397").Value = 90
$ws.Range("L397").Value = 100
$ws.Range("M397").Value = 95
$ws.Range("N397").Value = "$/unidad"
$ws.Range("O397").Value = "Región Metropolitana"
$ws.Range("P397").Value = 95
$ws.Range("Q397").Value = 1
$ws.Range("R397").Value = "Hortaliza"

# --- New row 398 ---
$ws.Range("A398").Value = 9
$ws.Range("B398").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C398").Value = "Metropolitana"
$ws.Range("D398").Value2 = 44610
$ws.Range("E398").Value = 13
$ws.Range("F398").Value = 100114014
$ws.Range("G398").Value = "Betarraga"
$ws.Range("H398").Value = "Sin especificar"
$ws.Range("I398").Value = "Segunda"
$ws.Range("J398").Value = 3400
$ws.Range("K398").Value = 70
$ws.Range("L398").Value = 70
$ws.Range("M398").Value = 70
$ws.Range("N398").Value = "$/unidad"
$ws.Range("O398").Value = "Región Metropolitana"
$ws.Range("P398").Value = 70
$ws.Range("Q398").Value = 1
$ws.Range("R398").Value = "Hortaliza"
